$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '28.028.14'
$ws.Cells.Item(2, 5).Value = '  +3.54%  '

$ws.Cells.Item(3, 4).Value = '1.804.16'
$ws.Cells.Item(3, 5).Value = '  +4.26%  '

$c = $ws.Cells.Item(4, 4)
$c.NumberFormat = "@"
$c.Value = '0.9981'
$ws.Cells.Item(4, 5).Value = '  -0.29%  '

$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = "@"
$c.Value = '315.57'
$ws.Cells.Item(5, 5).Value = '  +1.66%  '

$c = $ws.Cells.Item(6, 4)
$c.NumberFormat = "@"
$c.Value = '0.9986'
$ws.Cells.Item(6, 5).Value = '  -0.23%  '

$c = $ws.Cells.Item(7, 4)
$c.NumberFormat = "@"
$c.Value = '0.5439'
$ws.Cells.Item(7, 5).Value = '  +11.48%  '

$c = $ws.Cells.Item(8, 4)
$c.NumberFormat = "@"
$c.Value = '0.3795'
$ws.Cells.Item(8, 5).Value = '  +8.09%  '

$ws.Cells.Item(9, 5).Value = '  -0.80%  '

$ws.Cells.Item(10, 5).Value = '  +3.91%  '

$c = $ws.Cells.Item(11, 4)
$c.NumberFormat = "@"
$c.Value = '1.131'
$ws.Cells.Item(11, 5).Value = '  +7.73%  '

$c = $ws.Cells.Item(12, 4)
$c.NumberFormat = "@"
$c.Value = '0.9988'
$ws.Cells.Item(12, 5).Value = '  -0.27%  '

$c = $ws.Cells.Item(13, 4)
$c.NumberFormat = "@"
$c.Value = '21.13'
$ws.Cells.Item(13, 5).Value = '  +5.61%  '

$c = $ws.Cells.Item(14, 4)
$c.NumberFormat = "@"
$c.Value = '6.233'
$ws.Cells.Item(14, 5).Value = '  +5.84%  '

$ws.Cells.Item(15, 4).Value = '1.799.51'
$ws.Cells.Item(15, 5).Value = '  +3.92%  '

$c = $ws.Cells.Item(16, 4)
$c.NumberFormat = "@"
$c.Value = '7.143'
$ws.Cells.Item(16, 5).Value = '  +3.53%  '

$c = $ws.Cells.Item(17, 4)
$c.NumberFormat = "@"
$c.Value = '91.85'
$ws.Cells.Item(17, 5).Value = '  +5.15%  '

$c = $ws.Cells.Item(18, 4)
$c.NumberFormat = "@"
$c.Value = '0.00001078'
$ws.Cells.Item(18, 5).Value = '  +3.64%  '

$c = $ws.Cells.Item(19, 4)
$c.NumberFormat = "@"
$c.Value = '0.06508'
$ws.Cells.Item(19, 5).Value = '  +1.43%  '

$c = $ws.Cells.Item(20, 4)
$c.NumberFormat = "@"
$c.Value = '0.9986'
$ws.Cells.Item(20, 5).Value = '  -0.20%  '

$ws.Cells.Item(21, 5).Value = '  +3.25%  '

$c = $ws.Cells.Item(22, 4)
$c.NumberFormat = "@"
$c.Value = '5.980'
$ws.Cells.Item(22, 5).Value = '  +5.00%  '

$ws.Cells.Item(23, 4).Value = '28.025.67'
$ws.Cells.Item(23, 5).Value = '  +3.37%  '

$c = $ws.Cells.Item(24, 4)
$c.NumberFormat = "@"
$c.Value = '11.25'
$ws.Cells.Item(24, 5).Value = '  +2.85%  '

$c = $ws.Cells.Item(25, 4)
$c.NumberFormat = "@"
$c.Value = '2.098'
$ws.Cells.Item(25, 5).Value = '  +0.89%  '

$c = $ws.Cells.Item(26, 4)
$c.NumberFormat = "@"
$c.Value = '156.36'
$ws.Cells.Item(26, 5).Value = '  +1.69%  '

$c = $ws.Cells.Item(27, 4)
$c.NumberFormat = "@"
$c.Value = '20.57'
$ws.Cells.Item(27, 5).Value = '  +2.83%  '

$c = $ws.Cells.Item(28, 4)
$c.NumberFormat = "@"
$c.Value = '2.395'
$ws.Cells.Item(28, 5).Value = '  +14.59%  '

$ws.Cells.Item(29, 4).Value = '2.007.09'
$ws.Cells.Item(29, 5).Value = '  +4.01%  '

$c = $ws.Cells.Item(30, 4)
$c.NumberFormat = "@"
$c.Value = '122.91'
$ws.Cells.Item(30, 5).Value = '  +1.10%  '

$c = $ws.Cells.Item(31, 4)
$c.NumberFormat = "@"
$c.Value = '1.146'
$ws.Cells.Item(31, 5).Value = '  +9.30%  '

$c = $ws.Cells.Item(32, 4)
$c.NumberFormat = "@"
$c.Value = '0.1034'
$ws.Cells.Item(32, 5).Value = '  +10.88%  '

$c = $ws.Cells.Item(33, 4)
$c.NumberFormat = "@"
$c.Value = '5.758'
$ws.Cells.Item(33, 5).Value = '  +7.00%  '

$c = $ws.Cells.Item(34, 4)
$c.NumberFormat = "@"
$c.Value = '3.604'
$ws.Cells.Item(34, 5).Value = '  -1.03%  '

$c = $ws.Cells.Item(35, 4)
$c.NumberFormat = "@"
$c.Value = '0.02299'
$ws.Cells.Item(35, 5).Value = '  +4.99%  '

$ws.Cells.Item(36, 2).Value = 'Algorand'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$c = $ws.Cells.Item(36, 4)
$c.NumberFormat = "@"
$c.Value = '0.2126'
$ws.Cells.Item(36, 5).Value = '  +6.38%  '

$ws.Cells.Item(37, 2).Value = 'FraxShare'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$c = $ws.Cells.Item(37, 4)
$c.NumberFormat = "@"
$c.Value = '8.673'
$ws.Cells.Item(37, 5).Value = '  +15.47%  '

$ws.Cells.Item(38, 2).Value = 'Hedera'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$c = $ws.Cells.Item(38, 4)
$c.NumberFormat = "@"
$c.Value = '0.06062'
$ws.Cells.Item(38, 5).Value = '  +1.55%  '

$ws.Cells.Item(39, 2).Value = 'InternetComputer(DFINITY)'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$c = $ws.Cells.Item(39, 4)
$c.NumberFormat = "@"
$c.Value = '5.030'
$ws.Cells.Item(39, 5).Value = '  +5.36%  '

$c = $ws.Cells.Item(40, 4)
$c.NumberFormat = "@"
$c.Value = '11.52'
$ws.Cells.Item(40, 5).Value = '  +4.81%  '

$c = $ws.Cells.Item(41, 4)
$c.NumberFormat = "@"
$c.Value = '0.6286'
$ws.Cells.Item(41, 5).Value = '  +4.57%  '

$ws.Cells.Item(42, 2).Value = 'Frax'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$c = $ws.Cells.Item(42, 4)
$c.NumberFormat = "@"
$c.Value = '0.9987'
$ws.Cells.Item(42, 5).Value = '  -0.11%  '

$ws.Cells.Item(43, 2).Value = 'WEMIXTOKEN'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$c = $ws.Cells.Item(43, 4)
$c.NumberFormat = "@"
$c.Value = '1.399'
$ws.Cells.Item(43, 5).Value = '  -3.28%  '

$c = $ws.Cells.Item(44, 4)
$c.NumberFormat = "@"
$c.Value = '1.152'
$ws.Cells.Item(44, 5).Value = '  +5.09%  '

$c = $ws.Cells.Item(45, 4)
$c.NumberFormat = "@"
$c.Value = '13.41'
$ws.Cells.Item(45, 5).Value = '  +4.63%  '

$c = $ws.Cells.Item(46, 4)
$c.NumberFormat = "@"
$c.Value = '0.5929'
$ws.Cells.Item(46, 5).Value = '  +4.64%  '

$c = $ws.Cells.Item(47, 4)
$c.NumberFormat = "@"
$c.Value = '3.663'
$ws.Cells.Item(47, 5).Value = '  +2.22%  '

$c = $ws.Cells.Item(48, 4)
$c.NumberFormat = "@"
$c.Value = '122.24'
$ws.Cells.Item(48, 5).Value = '  +2.77%  '

$c = $ws.Cells.Item(49, 4)
$c.NumberFormat = "@"
$c.Value = '1.924'
$ws.Cells.Item(49, 5).Value = '  +3.88%  '

$ws.Cells.Item(50, 5).Value = '  +2.81%  '

$c = $ws.Cells.Item(51, 4)
$c.NumberFormat = "@"
$c.Value = '0.06787'
$ws.Cells.Item(51, 5).Value = '  +2.14%  '
